$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-10: shift each row's data to match the row that used to follow it
# (the original row 6 entry - "Processos Decisórios" - is removed, and
# everything below slides up by one).

# Row 6 (was row 7's data)
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "26/12/2025 06:11"
$ws.Range("C6").Value = 436
$ws.Range("E6").Value = "Logística"
$ws.Range("F6").Value = "Errando mais uma questão por confundir incorreta com correta"

# Row 7 (was row 8's data)
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "26/12/2025 06:21"
$ws.Range("C7").Value = 554
$ws.Range("E7").Value = "Engenharia de Métodos e Processos"
$ws.Range("F7").Value = "Investigar como é que aumento da eficiencia pode aumentar o tempo gasto ao inves de diminuir"

# Row 8 (was row 9's data)
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "27/12/2025 00:51"
$ws.Range("C8").Value = 542
$ws.Range("E8").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("F8").Value = "Estudar TPM"

# Row 9 (was row 10's data)
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "27/12/2025 00:52"
$ws.Range("F9").Value = "Estudar esse tal de 8 S"

# Row 10 (was row 11's data)
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "27/12/2025 05:40"
$ws.Range("C10").Value = 937
$ws.Range("D10").Value = "Inglês"
$ws.Range("E10").Value = "Pronouns"
$ws.Range("F10").Value = "Retirar os números de linhas no meio do testo"

# Row 11 now holds what used to be row 17's data (questao_id becomes text "62")
$ws.Range("A11").Value = 16
$ws.Range("B11").Value = "27/12/2025 06:38"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "62"
$ws.Range("D11").Value = "Português"
$ws.Range("E11").Value = "Conjunção"
$ws.Range("F11").Value = "Investigar o valor semântico de concessão e conformidade, marquei conformidade, letra D, nessa questão"

# Row 12 becomes a brand-new note (replacing the old "yield" note entirely)
$ws.Range("A12").Value = 21
$ws.Range("B12").Value = "27/12/2025 10:14"
$ws.Range("C12").Value = 121
$ws.Range("D12").Value = "Português"
$ws.Range("E12").Value = "Verbos Traiçoeiros"
$ws.Range("F12").Value = "Estudar mais esse assunto"

# Remove old rows 13-20, which are no longer part of the notebook
$ws.Range("A13:F20").ClearContents()
